# Natmi LR-pair result following Dr Hou advice:
# the sending/target cluster set grows from {FAPs, sCs} to {ECs, FAPs, sCs},
# so the 2x2 cluster-pair table (rows 2-5) becomes a full 3x3 table
# (rows 2-10) for the Fgf18 -> Fgfr2 ligand-receptor pair, with every
# metric column (E:T) recomputed for the new cluster composition.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 20
$dataRows = 9

$flat = @(
"ECs","Fgf18","Fgfr2","ECs",1,0.3333333333333333,0.4765323333333333,1.429597,0.04733670697480491,0.04733670697480491,
2,0.6666666666666666,0.09433999999999999,0.28302,0.05191071108246543,0.05191071108246543,0.04495606032666667,0.40460454294,0.002457282119364424,0.002457282119364424,
"ECs","Fgf18","Fgfr2","FAPs",1,0.3333333333333333,0.4765323333333333,1.429597,0.04733670697480491,0.04733670697480491,
3,1,0.9431116666666667,2.829335,0.5189484550226392,0.5189484550226391,0.4494232031105556,4.044808827995,0.0245653109504344,0.02456531095043439,
"ECs","Fgf18","Fgfr2","sCs",1,0.3333333333333333,0.4765323333333333,1.429597,0.04733670697480491,0.04733670697480491,
3,1,0.7798996666666667,2.339699,0.4291408338948954,0.4291408338948954,0.3716474079225556,3.344826671303,0.02031411390500609,0.02031411390500609,
"FAPs","Fgf18","Fgfr2","ECs",3,1,8.770835333333332,26.312506,0.8712576945074703,0.8712576945074703,
2,0.6666666666666666,0.09433999999999999,0.28302,0.05191071108246543,0.05191071108246543,0.8274406053466665,7.446965448119999,0.04522760645795221,0.04522760645795221,
"FAPs","Fgf18","Fgfr2","FAPs",3,1,8.770835333333332,26.312506,0.8712576945074703,0.8712576945074703,
3,1,0.9431116666666667,2.829335,0.5189484550226392,0.5189484550226391,8.271877129278888,74.44689416351,0.4521378344912382,0.4521378344912382,
"FAPs","Fgf18","Fgfr2","sCs",3,1,8.770835333333332,26.312506,0.8712576945074703,0.8712576945074703,
3,1,0.7798996666666667,2.339699,0.4291408338948954,0.4291408338948954,6.840371552854888,61.563343975694,0.3738922535582798,0.3738922535582798,
"sCs","Fgf18","Fgfr2","ECs",3,1,0.8194993333333332,2.458498,0.08140559851772487,0.08140559851772486,
2,0.6666666666666666,0.09433999999999999,0.28302,0.05191071108246543,0.05191071108246543,0.07731156710666666,0.6958041039599999,0.004225822505148791,0.004225822505148791,
"sCs","Fgf18","Fgfr2","FAPs",3,1,0.8194993333333332,2.458498,0.08140559851772487,0.08140559851772486,
3,1,0.9431116666666667,2.829335,0.5189484550226392,0.5189484550226391,0.7728793820922222,6.955914438829999,0.04224530958096657,0.04224530958096655,
"sCs","Fgf18","Fgfr2","sCs",3,1,0.8194993333333332,2.458498,0.08140559851772487,0.08140559851772486,
3,1,0.7798996666666667,2.339699,0.4291408338948954,0.4291408338948954,0.6391272569002221,5.752145312101999,0.03493446643160951,0.0349344664316095
)

for ($r = 0; $r -lt $dataRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $idx = $r*$numCols + $c
        $ws.Cells.Item($r+2, $c+1).Value = $flat[$idx]
    }
}
